$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reset sheet: remove all values but keep/expand used range appropriately ---
$ws.Cells.ClearContents()
$ws.Range("X1:AG19").Delete()

# --- Row 1 header: sequential numbers 0..21 across B1:W1 ---
$row1 = New-Object 'object[,]' 1,22
for ($i = 0; $i -lt 22; $i++) { $row1[0,$i] = $i }
$ws.Range("B1:W1").Value = $row1

# --- Column B (B2:B29): HKL-family labels. Written top-to-bottom so that the
#     shared-string table gets these 28 strings in this exact order (indices 0-27). ---
$colB = New-Object 'object[,]' 28,1
$colB[0,0] = "HKL"
$colB[1,0] = "Spiral5"
$colB[2,0] = "RotRing OmegaMax-90"
$colB[3,0] = "Equal Angle"
$colB[4,0] = "Tilt Rotate"
$colB[5,0] = "CLR"
$colB[6,0] = "Rizzie Hex"
$colB[7,0] = "Thomas Hex"
$colB[8,0] = "Tilt Rotate_Partial"
$colB[9,0] = "RotRing OmegaMax-60"
$colB[10,0] = "Equal Angle_Partial"
$colB[11,0] = "Rizzie Hex_Partial"
$colB[12,0] = "ND Single"
$colB[13,0] = "RD Single"
$colB[14,0] = "TD Single"
$colB[15,0] = "Morris Single"
$colB[16,0] = "Ring Perpendicular to ND"
$colB[17,0] = "Ring Perpendicular to RD"
$colB[18,0] = "Ring Perpendicular to TD"
$colB[19,0] = "OffsetFTD"
$colB[20,0] = "OffsetATD"
$colB[21,0] = "OffsetF45"
$colB[22,0] = "OffsetA45"
$colB[23,0] = "OffsetFRD"
$colB[24,0] = "OffsetARD"
$colB[25,0] = "Gaussian Quadrature"
$colB[26,0] = "Michael-CCHex"
$colB[27,0] = "Michael-SNHex"
$ws.Range("B2:B29").Value = $colB

# --- Row 2 (C2:W2): HKL-index / pair labels, in order, continuing the shared-string
#     table with indices 28-48. ---
$row2 = New-Object 'object[,]' 1,21
$row2[0,0] = "[2, 0, 0]"
$row2[0,1] = "[2, 2, 0]"
$row2[0,2] = "[3, 3, 3]"
$row2[0,3] = "[4, 2, 0]"
$row2[0,4] = "[4, 0, 0]"
$row2[0,5] = "[4, 2, 2]"
$row2[0,6] = "[5, 1, 1]"
$row2[0,7] = "[1, 1, 1]"
$row2[0,8] = "[2, 2, 2]"
$row2[0,9] = "[3, 3, 1]"
$row2[0,10] = "[3, 1, 1]"
$row2[0,11] = "1Pair-A"
$row2[0,12] = "1Pair-B"
$row2[0,13] = "2Pairs-A"
$row2[0,14] = "2Pairs-B"
$row2[0,15] = "3Pairs-A"
$row2[0,16] = "3Pairs-B"
$row2[0,17] = "3Pairs-C"
$row2[0,18] = "4Pairs"
$row2[0,19] = "5A4F"
$row2[0,20] = "MaxUnique"
$ws.Range("C2:W2").Value = $row2

# --- Column A (A2:A29): sequential numbers 0..27 ---
$colA = New-Object 'object[,]' 28,1
for ($i = 0; $i -lt 28; $i++) { $colA[$i,0] = $i }
$ws.Range("A2:A29").Value = $colA

# --- Body C3:W29 (27 rows x 21 cols) filled with 1 ---
$body = New-Object 'object[,]' 27,21
for ($r = 0; $r -lt 27; $r++) { for ($c = 0; $c -lt 21; $c++) { $body[$r,$c] = 1 } } 
$ws.Range("C3:W29").Value = $body

# --- Styling: re-apply the bold/border/center header style to the cells that need it.
#     B1:W1 and A2:A19 already retained style "1" from ClearContents(); we only need to
#     extend that same style down to the newly added A20:A29. ---
$ws.Range("A2").Copy()
$ws.Range("A20:A29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host "done"
